$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Insert the three new columns ------------------------------------------
# Original layout:  A B C D E(especialidade) F(turno) G(suplentes)
# Target layout:     A B C D E(matr_apropriador) F(nome_setor) G(nome_planilha) H(turno) I(suplentes) J(flags)
#
# Insert a blank column in front of the old "especialidade" column (E) so it
# shifts right to F (keeping its data, e.g. "SOLDA", intact).
$ws.Columns.Item(5).Insert()
# Insert a second blank column right after the shifted "especialidade" column
# (now F) so the following turno/suplentes columns move to H/I.
$ws.Columns.Item(7).Insert()
# Insert a new trailing blank column after suplentes (I) for "flags".
$ws.Columns.Item(10).Insert()

# The inserted columns carry over formatting from their neighbours; make sure
# the new data cells (rows 2-6) are completely empty, matching the target.
$ws.Range("E2:E6").Clear()
$ws.Range("G2:G6").Clear()
$ws.Range("J2:J6").Clear()

# --- Header row text ---------------------------------------------------------
$ws.Range("E1").Value = "matr_apropriador"
$ws.Range("F1").Value = "nome_setor"
$ws.Range("G1").Value = "nome_planilha"
$ws.Range("J1").Value = "flags"

# --- Header row formatting (yellow for new blank columns, green for nome_setor)
$ws.Range("E1").Interior.Color = 65535
$ws.Range("E1").Borders.LineStyle = 1

$ws.Range("F1").Interior.Color = 5296274
$ws.Range("F1").Borders.LineStyle = 1

$ws.Range("G1").Interior.Color = 65535
$ws.Range("G1").Borders.LineStyle = 1

$ws.Range("J1").Interior.Color = 65535
$ws.Range("J1").Borders.LineStyle = 1

# --- Column widths for the brand-new columns (E, G); F/H/I keep the widths
# that were carried over automatically from their shifted source columns.
$ws.Columns.Item(5).ColumnWidth = 17
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666

# --- Defined name / filter database now spans the extra columns ------------
$wb.Names.Item(1).RefersTo = "=Plan1!`$A`$1:`$I`$6"

Write-Host "done"
